$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the IG generation timestamp ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(8, 2).Value = "2025-07-24T13:17:05+00:00"

# --- Elements sheet: append a new row describing ExerciceProfessionnel.professionnel ---
$ws = $wb.Worksheets.Item("Elements")

$lastRow = 28
$newRow  = 29

# Seed the new row as a full copy (values + styling) of the last existing
# leaf row. Every row in this table shares one cell style ("s=2"), and most
# columns are blank/"0"/"1" placeholders, so cloning the row gives us the
# right style + the right "blank vs placeholder" pattern for free.
$ws.Range("A$lastRow`:AJ$lastRow").Copy($ws.Range("A$newRow`:AJ$newRow"))

# Min/Base Min (col F / col AG) need to read "1" instead of the copied "0".
# Pull the text "1" from this row's own Max/Base Max cell (col G / col AH)
# instead of typing a bare "1" literal, so it stays a text cell (matching
# the sheet's convention) rather than being auto-coerced to a number.
$ws.Cells.Item($lastRow, 7).Copy($ws.Cells.Item($newRow, 6))
$ws.Cells.Item($lastRow, 7).Copy($ws.Cells.Item($newRow, 33))

# Now overwrite the element-specific text for the new property.
$id = "ExerciceProfessionnel.professionnel"
$type = "Reference(https://interop.esante.gouv.fr/ig/fhir/mos/StructureDefinition/Professionnel)`n"
$short = "Lien vers la classe Professionnel."

$ws.Cells.Item($newRow, 1).Value = $id      # ID
$ws.Cells.Item($newRow, 2).Value = $id      # Path
$ws.Cells.Item($newRow, 11).Value = $type   # Type(s)
$ws.Cells.Item($newRow, 12).Value = $short  # Short
$ws.Cells.Item($newRow, 13).Value = $short  # Definition
$ws.Cells.Item($newRow, 32).Value = $id     # Base Path

# Widen column K ("Type(s)") so the longer Reference(...) text keeps fitting
# ("bestFit" semantics from the original file, approximated via ColumnWidth).
$ws.Columns.Item(11).ColumnWidth = 68.0

# Writing the long, multi-line Type(s)/Short/Definition text triggers an
# auto row-height bump (wrap text). The source workbook has no explicit row
# heights anywhere, so re-autofit the row to drop back to the implicit
# default height instead of leaving a stray customHeight behind.
$ws.Rows.Item($newRow).AutoFit()
